$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "John"
$ws.Range("A3").Value = "James"
$ws.Range("A4").Value = "Scott "
$ws.Range("A5").Value = "John"

$ws.Range("B2").Value = "Doe"
$ws.Range("B3").Value = "Smith"
$ws.Range("B4").Value = "Patterson"
$ws.Range("B5").Value = "Ditto"

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "0000001"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "0000003"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "0000004"

$ws.Range("B10").Select()
